# M1 C3 P2 e M2 C2 P1
#
# 1. Update the cached "today" date field text (datetimeFigureOut) on the
#    slide master and every slide layout: "8/8/2015" -> "15/8/2015".
# 2. Merge the two runs "Fluxo de " + "Código" on slide 2 (Sumário) into a
#    single run "Fluxo de Código".
# 3. Remove the "Exercícios com if" slide (old position 44) and the
#    "Exercícios com switch" slide (old position 47) - the "Tipo enum" and
#    "Operador Condicional switch" slides that used to sit between them
#    naturally shift up to take positions 44 and 45.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $t = $sh.TextFrame.TextRange.Text
            if ($t -eq "8/8/2015") {
                $sh.TextFrame.TextRange.Text = "15/8/2015"
            }
        }
    }
}

# -- 1. date placeholder caches on master + all layouts --------------------
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# -- 2. merge "Fluxo de " + "Código" into one run on slide 2 ---------------
$summarySlide = $p.Slides.Item(2)
$contentShape = $summarySlide.Shapes.Item(2)
$firstPara = $contentShape.TextFrame.TextRange.Paragraphs(1, 1)
$firstPara.Text = "Fluxo de Código"

# -- 3. drop the two standalone exercise slides -----------------------------
$exerciseIfSlide = $p.Slides.Item(44)
$exerciseIfSlide.Delete()

$exerciseSwitchSlide = $p.Slides.Item(46)
$exerciseSwitchSlide.Delete()
